$wb = $excel.ActiveWorkbook

# --- Sheet2: LoginPageTest ---
# New string "admin12" must be interned into sharedStrings before "N" (see below)
# so that shared-string indices come out in the same order as the target file.
$ws2 = $wb.Worksheets.Item("LoginPageTest")
$ws2.Range("A2:C2").Copy()
$ws2.Range("A3:C3").PasteSpecial(-4122)
$ws2.Range("A3").Value = "Admin"
$ws2.Range("B3").Value = "admin12"
$ws2.Range("C3").Value = "Y"

# --- Sheet1: TestSuite ---
$ws1 = $wb.Worksheets.Item("TestSuite")
$ws1.Activate()
$ws1.Range("B3").Value = "N"
$ws1.Range("F3").Select()

# --- Sheet3: AdminPageTest ---
$ws3 = $wb.Worksheets.Item("AdminPageTest")
$ws3.Activate()
$ws3.Range("E2").Select()

# Re-activate LoginPageTest last so it remains the selected tab, and set its
# own selection, matching the original workbook state.
$ws2.Activate()
$ws2.Range("A3").Select()
